# Update mods data [2026-01-24 15:09:22]
# Append a new daily snapshot row (row 75) to the ModCounts sheet,
# mirroring the existing rows: Date (text), Game (text), ModCount (number).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 75

# Write the date as literal text (leading apostrophe forces text entry so
# Excel doesn't reinterpret "2026/01/24" as a date serial number), then the
# remaining columns.
$ws.Cells.Item($newRow, 1).Formula = "'2026/01/24"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1156

# Match the formatting (centered alignment style) used by the rest of the
# data rows by copying the previous row's format onto the new one.
$ws.Range("A74:C74").Copy()
$ws.Range("A75:C75").PasteSpecial(-4122)
